# Auto-generated: apply 2022-11-07 crime data updates to violent-crime-full-year.xlsx
# 143 cell updates across 42 worksheets (column I = year 2022 totals).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 6264
$ws.Range("I3").Value = 6535
$ws.Range("I4").Value = 1501
$ws.Range("I5").Value = 605
$ws.Range("I6").Value = 7410
$ws.Range("I7").Value = 22315

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I3").Value = 63
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 44
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 127

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 209
$ws.Range("I3").Value = 235
$ws.Range("I6").Value = 207
$ws.Range("I7").Value = 705

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 320
$ws.Range("I5").Value = 25
$ws.Range("I7").Value = 862

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 175
$ws.Range("I7").Value = 705
$ws.Range("I8").Value = 1345
$ws.Range("I10").Value = 158
$ws.Range("I11").Value = 333
$ws.Range("I12").Value = 53
$ws.Range("I14").Value = 127
$ws.Range("I18").Value = 168
$ws.Range("I19").Value = 624
$ws.Range("I21").Value = 100
$ws.Range("I22").Value = 63
$ws.Range("I27").Value = 195
$ws.Range("I29").Value = 1359
$ws.Range("I33").Value = 1015
$ws.Range("I34").Value = 103
$ws.Range("I36").Value = 303
$ws.Range("I37").Value = 705
$ws.Range("I42").Value = 783
$ws.Range("I43").Value = 191
$ws.Range("I47").Value = 162
$ws.Range("I51").Value = 267
$ws.Range("I52").Value = 480
$ws.Range("I53").Value = 243
$ws.Range("I54").Value = 457
$ws.Range("I61").Value = 23
$ws.Range("I62").Value = 7
$ws.Range("I63").Value = 70
$ws.Range("I66").Value = 63
$ws.Range("I67").Value = 862
$ws.Range("I70").Value = 35
$ws.Range("I76").Value = 322
$ws.Range("I78").Value = 303
$ws.Range("I79").Value = 630
$ws.Range("I83").Value = 485
$ws.Range("I85").Value = 1006
$ws.Range("I86").Value = 138
$ws.Range("I88").Value = 205
$ws.Range("I89").Value = 260
$ws.Range("I90").Value = 281
$ws.Range("I91").Value = 234
$ws.Range("I101").Value = 22315

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 177
$ws.Range("I7").Value = 485

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 230
$ws.Range("I3").Value = 377
$ws.Range("I6").Value = 321
$ws.Range("I7").Value = 1015

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 218
$ws.Range("I7").Value = 457

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 400
$ws.Range("I6").Value = 377
$ws.Range("I7").Value = 1359

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 208
$ws.Range("I3").Value = 187
$ws.Range("I7").Value = 624

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 66
$ws.Range("I6").Value = 145
$ws.Range("I7").Value = 322

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 286
$ws.Range("I6").Value = 253
$ws.Range("I7").Value = 1006

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I6").Value = 266
$ws.Range("I7").Value = 783

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 71
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I4").Value = 40
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 234

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 184
$ws.Range("I3").Value = 199
$ws.Range("I6").Value = 188
$ws.Range("I7").Value = 630

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I6").Value = 129
$ws.Range("I7").Value = 480

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I3").Value = 28
$ws.Range("I6").Value = 22
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 162

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I3").Value = 67
$ws.Range("I7").Value = 333

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 64
$ws.Range("I7").Value = 175

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 61
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 402
$ws.Range("I3").Value = 382
$ws.Range("I6").Value = 436
$ws.Range("I7").Value = 1345

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 64
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 97
$ws.Range("I7").Value = 281

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 73
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 267

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 55
$ws.Range("I3").Value = 53
$ws.Range("I6").Value = 110
$ws.Range("I7").Value = 243

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 233
$ws.Range("I7").Value = 705

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("I2").Value = 2
$ws.Range("I6").Value = 7

Write-Output "Applied 2022-11-07 crime data updates."